# 20250515 without elite inheritance
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-16 (Run=1 Solution_Index 7-8, and all of Run=2 rows)
$ws.Range("A8:G16").EntireRow.Delete() | Out-Null

# Update rows 2-7 with the new values from the re-run without elite inheritance
$data = @(
    @(1, 1, 0.2522901708343649, 0.8817639482939641, "{641, 2631, 2183, 2601, 2986, 107, 2351, 1622, 2201, 2973}", 2.352813243865967, 0.5819225762028103),
    @(1, 2, 0.254518445159693,  0.8273204540573542, "{2273, 2183, 1320, 107, 2384, 3826, 2706, 2550, 3101, 2462}", 2.352813243865967, 0.5819225762028103),
    @(1, 3, 0.1631344392176281, 0.9127066249820019, "{3907, 2244, 3366, 1001, 1515, 2348, 2292, 1049, 3101, 414}", 2.352813243865967, 0.5819225762028103),
    @(2, 1, 0.2778905669720228, 0.8756980904083504, "{3201, 2694, 1864, 107, 3051, 109, 2353, 1076, 2038, 1912}",  1.589942455291748, 0.6327568357026606),
    @(2, 2, 0.2393414211438475, 0.8962477314085928, "{1447, 136, 1450, 107, 2090, 1486, 1616, 2706, 122, 2619}",   1.589942455291748, 0.6327568357026606),
    @(2, 3, 0.1948502104481307, 0.897609454607404,  "{0, 1536, 2500, 1126, 1577, 363, 3116, 1197, 239, 1758}",     1.589942455291748, 0.6327568357026606)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
